$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

# Row 19
$ws.Range("H19").Value = 709.86957
$ws.Range("I19").Value = 584.2222
$ws.Range("K19").Value = 584.2222
$ws.Range("M19").Value = -409.2222

# Row 53
$ws.Range("H53").Value = 104.625
$ws.Range("I53").Value = 78.85714
$ws.Range("J53").Value = 124.666664
$ws.Range("K53").Value = 78.85714
$ws.Range("L53").Value = 124.666664
$ws.Range("M53").Value = 558.14286
$ws.Range("N53").Value = -1398.666664

# Row 80
$ws.Range("H80").Value = 6846.7646
$ws.Range("I80").Value = 349.5
$ws.Range("J80").Value = 12622.111
$ws.Range("K80").Value = 1048.5
$ws.Range("L80").Value = 37866.333
$ws.Range("M80").Value = -50.5
$ws.Range("N80").Value = -39862.333

# Row 83
$ws.Range("H83").Value = 6846.7646
$ws.Range("I83").Value = 349.5
$ws.Range("J83").Value = 12622.111
$ws.Range("K83").Value = 3145.5
$ws.Range("L83").Value = 113598.999
$ws.Range("M83").Value = 1846.5
$ws.Range("N83").Value = -123582.999

# Row 86
$ws.Range("H86").Value = 2421.8462
$ws.Range("I86").Value = 2220.4443
$ws.Range("K86").Value = 2220.4443
$ws.Range("M86").Value = -1097.4443

# Row 89
$ws.Range("H89").Value = 2421.8462
$ws.Range("I89").Value = 2220.4443
$ws.Range("K89").Value = 11102.2215
$ws.Range("M89").Value = -5486.2215

# Row 125
$ws.Range("H125").Value = 3012.0908
$ws.Range("I125").Value = 3223.75
$ws.Range("K125").Value = 29013.75
$ws.Range("M125").Value = -26553.75


$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4745.778
$ws.Range("I61").Value = 5428
$ws.Range("K61").Value = 5428
$ws.Range("M61").Value = -5216

# Row 102
$ws.Range("H102").Value = 3290
$ws.Range("I102").Value = 3158.5715
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 3158.5715
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -1536.5715
$ws.Range("N102").Value = -6994

# Row 123
$ws.Range("H123").Value = 25628
$ws.Range("J123").Value = 25628
$ws.Range("L123").Value = 25628
$ws.Range("N123").Value = -35428

# Row 136
$ws.Range("H136").Value = 4745.778
$ws.Range("I136").Value = 5428
$ws.Range("K136").Value = 16284
$ws.Range("M136").Value = -13734


$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 4999.75
$ws.Range("I107").Value = 5000
$ws.Range("J107").Value = 4999
$ws.Range("K107").Value = 5000
$ws.Range("L107").Value = 4999
$ws.Range("M107").Value = -3080
$ws.Range("N107").Value = -8839

# Row 134
$ws.Range("H134").Value = 4510.4375
$ws.Range("J134").Value = 4185.5
$ws.Range("L134").Value = 12556.5
$ws.Range("N134").Value = -17626.5


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3672.2
$ws.Range("I31").Value = 1840.25
$ws.Range("J31").Value = 11000
$ws.Range("K31").Value = 1840.25
$ws.Range("L31").Value = 11000
$ws.Range("M31").Value = -1545.25
$ws.Range("N31").Value = -11590

# Row 34
$ws.Range("H34").Value = 3672.2
$ws.Range("I34").Value = 1840.25
$ws.Range("J34").Value = 11000
$ws.Range("K34").Value = 1840.25
$ws.Range("L34").Value = 11000
$ws.Range("M34").Value = -1638.25
$ws.Range("N34").Value = -11404

# Row 62
$ws.Range("H62").Value = 74129.28999999999
$ws.Range("I62").Value = 102241
$ws.Range("K62").Value = 102241
$ws.Range("M62").Value = -101617

# Row 65
$ws.Range("H65").Value = 74129.28999999999
$ws.Range("I65").Value = 102241
$ws.Range("K65").Value = 511205
$ws.Range("M65").Value = -508085

# Row 107
$ws.Range("H107").Value = 383.41666
$ws.Range("I107").Value = 202.75
$ws.Range("J107").Value = 473.75
$ws.Range("K107").Value = 202.75
$ws.Range("L107").Value = 473.75
$ws.Range("M107").Value = 1717.25
$ws.Range("N107").Value = -4313.75

# Row 132
$ws.Range("H132").Value = 2395.6155
$ws.Range("I132").Value = 1827.0555
$ws.Range("J132").Value = 3674.875
$ws.Range("K132").Value = 5481.166499999999
$ws.Range("L132").Value = 11024.625
$ws.Range("M132").Value = -2951.166499999999
$ws.Range("N132").Value = -16084.625


$ws = $wb.Worksheets.Item("CUL")
# Row 100
$ws.Range("H100").Value = 26500
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 9000
$ws.Range("N100").Value = -10622

# Row 103
$ws.Range("H103").Value = 2553.889
$ws.Range("J103").Value = 2808.2666
$ws.Range("L103").Value = 8424.799800000001
$ws.Range("N103").Value = -10182.7998

# Row 116
$ws.Range("H116").Value = 800
$ws.Range("I116").Value = 800
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2400
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1042
$ws.Range("N116").ClearContents()

# Row 131
$ws.Range("H131").Value = 35715844
$ws.Range("J131").Value = 45456336
$ws.Range("L131").Value = 136369008
$ws.Range("N131").Value = -136379088


$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3650
$ws.Range("I80").Value = 3527.7778
$ws.Range("J80").Value = 3925
$ws.Range("K80").Value = 3527.7778
$ws.Range("L80").Value = 3925
$ws.Range("M80").Value = -2529.7778
$ws.Range("N80").Value = -5921

# Row 83
$ws.Range("H83").Value = 3650
$ws.Range("I83").Value = 3527.7778
$ws.Range("J83").Value = 3925
$ws.Range("K83").Value = 17638.889
$ws.Range("L83").Value = 19625
$ws.Range("M83").Value = -12646.889
$ws.Range("N83").Value = -29609

# Row 102
$ws.Range("H102").Value = 2371.95
$ws.Range("I102").Value = 2409.4167
$ws.Range("J102").Value = 2315.75
$ws.Range("K102").Value = 2409.4167
$ws.Range("L102").Value = 2315.75
$ws.Range("M102").Value = -787.4167000000002
$ws.Range("N102").Value = -5559.75

# Row 107
$ws.Range("H107").Value = 1434.8889
$ws.Range("I107").Value = 2458.6667
$ws.Range("J107").Value = 411.1111
$ws.Range("K107").Value = 2458.6667
$ws.Range("L107").Value = 411.1111
$ws.Range("M107").Value = -538.6667000000002
$ws.Range("N107").Value = -4251.1111

# Row 132
$ws.Range("H132").Value = 3748.0908
$ws.Range("I132").Value = 2709.889
$ws.Range("K132").Value = 8129.667
$ws.Range("M132").Value = -5599.667


$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 697.5
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 896.6667
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 896.6667
$ws.Range("M55").Value = 73
$ws.Range("N55").Value = -1242.6667


$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 37881150
$ws.Range("I122").Value = 59525360
$ws.Range("J122").Value = 3790
$ws.Range("K122").Value = 178576080
$ws.Range("L122").Value = 11370
$ws.Range("M122").Value = -178573630
$ws.Range("N122").Value = -16270

# Row 136
$ws.Range("H136").Value = 2870.75
$ws.Range("I136").Value = 2422.5715
$ws.Range("J136").Value = 3916.5
$ws.Range("K136").Value = 7267.7145
$ws.Range("L136").Value = 11749.5
$ws.Range("M136").Value = -4717.7145
$ws.Range("N136").Value = -16849.5

